$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers: B1 "blood" -> "sugarBlood", C1 "dofek" -> "pulse"
$ws.Range("B1").Value = "sugarBlood"
$ws.Range("C1").Value = "pulse"

# Delete column E entirely (the "check" column, all filled with 4s)
$ws.Columns("E").Delete()

# Append new rows 12-18
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 100
$ws.Range("C12").Value = 80
$ws.Range("D12").Value = 42

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 61
$ws.Range("C13").Value = 88
$ws.Range("D13").Value = 24

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 72
$ws.Range("C14").Value = 78
$ws.Range("D14").Value = 31

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 82
$ws.Range("C15").Value = 76
$ws.Range("D15").Value = 40

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 115
$ws.Range("C16").Value = 86
$ws.Range("D16").Value = 64

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 70
$ws.Range("C17").Value = 68
$ws.Range("D17").Value = 54

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 56
$ws.Range("C18").Value = 84
$ws.Range("D18").Value = 27

# Update selection to J7 (matches post-edit sheetView selection in target)
$ws.Range("J7").Select()
